# Scheduled market-data refresh: update currentAveragePrice* / LevePrice* /
# LeveProfit* columns (H:N) for a handful of leve rows across all eight
# crafting-job sheets, per upstream Universalis price pull.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 3799.8
$ws.Range("I74").Value = 3666.3333
$ws.Range("J74").Value = 4000
$ws.Range("K74").Value = 3666.3333
$ws.Range("L74").Value = 4000
$ws.Range("M74").Value = -2730.3333
$ws.Range("N74").Value = -5872
$ws.Range("H77").Value = 3799.8
$ws.Range("I77").Value = 3666.3333
$ws.Range("J77").Value = 4000
$ws.Range("K77").Value = 18331.6665
$ws.Range("L77").Value = 20000
$ws.Range("M77").Value = -13651.6665
$ws.Range("N77").Value = -29360
$ws.Range("H138").Value = 2668533.2
$ws.Range("I138").Value = 1029.2285
$ws.Range("J138").Value = 5002599
$ws.Range("K138").Value = 3087.6855
$ws.Range("L138").Value = 15007797
$ws.Range("M138").Value = 2052.3145
$ws.Range("N138").Value = -15018077
$ws.Range("H139").Value = 33984.8
$ws.Range("J139").Value = 37481.25
$ws.Range("L139").Value = 37481.25
$ws.Range("N139").Value = -47761.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1941
$ws.Range("I45").Value = 2046.5834
$ws.Range("J45").Value = 1800.2222
$ws.Range("K45").Value = 2046.5834
$ws.Range("L45").Value = 1800.2222
$ws.Range("M45").Value = -1669.5834
$ws.Range("N45").Value = -2554.2222

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3643.173
$ws.Range("I134").Value = 2896.15
$ws.Range("J134").Value = 6133.25
$ws.Range("K134").Value = 8688.450000000001
$ws.Range("L134").Value = 18399.75
$ws.Range("M134").Value = -6153.450000000001
$ws.Range("N134").Value = -23469.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 32259458
$ws.Range("I31").Value = 43479176
$ws.Range("J31").Value = 2763.875
$ws.Range("K31").Value = 43479176
$ws.Range("L31").Value = 2763.875
$ws.Range("M31").Value = -43478881
$ws.Range("N31").Value = -3353.875
$ws.Range("H34").Value = 32259458
$ws.Range("I34").Value = 43479176
$ws.Range("J34").Value = 2763.875
$ws.Range("K34").Value = 43479176
$ws.Range("L34").Value = 2763.875
$ws.Range("M34").Value = -43478974
$ws.Range("N34").Value = -3167.875
$ws.Range("H51").Value = 16000
$ws.Range("J51").Value = 16000
$ws.Range("L51").Value = 16000
$ws.Range("N51").Value = -17472
$ws.Range("H59").Value = 10009.583
$ws.Range("J59").Value = 10009.583
$ws.Range("L59").Value = 10009.583
$ws.Range("N59").Value = -12299.583
$ws.Range("H60").Value = 11103
$ws.Range("J60").Value = 11103
$ws.Range("L60").Value = 11103
$ws.Range("N60").Value = -12125
$ws.Range("H61").Value = 16000
$ws.Range("J61").Value = 16000
$ws.Range("L61").Value = 16000
$ws.Range("N61").Value = -16696
$ws.Range("H96").Value = 27000
$ws.Range("J96").Value = 27000
$ws.Range("L96").Value = 27000
$ws.Range("N96").Value = -32492
$ws.Range("H132").Value = 3130.8125
$ws.Range("I132").Value = 2469.2222
$ws.Range("J132").Value = 3981.4285
$ws.Range("K132").Value = 7407.6666
$ws.Range("L132").Value = 11944.2855
$ws.Range("M132").Value = -4877.6666
$ws.Range("N132").Value = -17004.2855

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H64").Value = 4771.615
$ws.Range("I64").Value = 3487.1428
$ws.Range("J64").Value = 6270.1665
$ws.Range("K64").Value = 10461.4284
$ws.Range("L64").Value = 18810.4995
$ws.Range("M64").Value = -10191.4284
$ws.Range("N64").Value = -19350.4995
$ws.Range("H67").Value = 4771.615
$ws.Range("I67").Value = 3487.1428
$ws.Range("J67").Value = 6270.1665
$ws.Range("K67").Value = 10461.4284
$ws.Range("L67").Value = 18810.4995
$ws.Range("M67").Value = -9525.428400000001
$ws.Range("N67").Value = -20682.4995
$ws.Range("H75").Value = 3215.4
$ws.Range("J75").Value = 5270.8
$ws.Range("L75").Value = 15812.4
$ws.Range("N75").Value = -17808.4
$ws.Range("H78").Value = 3215.4
$ws.Range("J78").Value = 5270.8
$ws.Range("L78").Value = 47437.2
$ws.Range("N78").Value = -57421.2
$ws.Range("H131").Value = 879.62
$ws.Range("I131").Value = 564.875
$ws.Range("J131").Value = 906.98914
$ws.Range("K131").Value = 1694.625
$ws.Range("L131").Value = 2720.96742
$ws.Range("M131").Value = 3345.375
$ws.Range("N131").Value = -12800.96742

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H63").Value = 24000
$ws.Range("J63").Value = 26000
$ws.Range("L63").Value = 26000
$ws.Range("N63").Value = -27372
$ws.Range("H66").Value = 24000
$ws.Range("J66").Value = 26000
$ws.Range("L66").Value = 78000
$ws.Range("N66").Value = -84864
$ws.Range("H126").Value = 2012
$ws.Range("I126").Value = 2017.6
$ws.Range("J126").Value = 1993.3334
$ws.Range("K126").Value = 6052.799999999999
$ws.Range("L126").Value = 5980.0002
$ws.Range("M126").Value = -3582.799999999999
$ws.Range("N126").Value = -10920.0002

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H18").Value = 60000
$ws.Range("I18").Value = 0
$ws.Range("J18").Value = 60000
$ws.Range("K18").Value = 0
$ws.Range("L18").Value = 60000
$ws.Range("M18").ClearContents()
$ws.Range("N18").Value = -60344
$ws.Range("H20").Value = 40002.5
$ws.Range("I20").Value = 40002.5
$ws.Range("J20").Value = 0
$ws.Range("K20").Value = 40002.5
$ws.Range("L20").Value = 0
$ws.Range("M20").Value = -39776.5
$ws.Range("N20").ClearContents()
$ws.Range("H22").Value = 44204.824
$ws.Range("I22").Value = 1000000
$ws.Range("J22").Value = 759.5909
$ws.Range("K22").Value = 1000000
$ws.Range("L22").Value = 759.5909
$ws.Range("M22").Value = -999705
$ws.Range("N22").Value = -1349.5909
$ws.Range("H27").Value = 44204.824
$ws.Range("I27").Value = 1000000
$ws.Range("J27").Value = 759.5909
$ws.Range("K27").Value = 1000000
$ws.Range("L27").Value = 759.5909
$ws.Range("M27").Value = -999893
$ws.Range("N27").Value = -973.5909
$ws.Range("H136").Value = 7937650.5
$ws.Range("I136").Value = 9804960
$ws.Range("J136").Value = 1584.625
$ws.Range("K136").Value = 29414880
$ws.Range("L136").Value = 4753.875
$ws.Range("M136").Value = -29412330
$ws.Range("N136").Value = -9853.875

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H64").Value = 27990
$ws.Range("J64").Value = 27990
$ws.Range("L64").Value = 27990
$ws.Range("N64").Value = -28486
$ws.Range("H67").Value = 27990
$ws.Range("J67").Value = 27990
$ws.Range("L67").Value = 27990
$ws.Range("N67").Value = -29706
$ws.Range("H137").Value = 86286
$ws.Range("J137").Value = 86286
$ws.Range("L137").Value = 86286
$ws.Range("N137").Value = -96486
$ws.Range("H140").Value = 34163
$ws.Range("J140").Value = 34163
$ws.Range("L140").Value = 34163
$ws.Range("N140").Value = -44523
$ws.Range("H141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("L141").Value = 0
